# Workbook: 化学纤维.xlsx
# Changes applied:
#  1. Within each 4-row year block (A/B/C/D sub-rows), the "B" and "C"
#     sub-rows swap places content-wise (labels + data), while the row
#     position in the sheet (the row number) stays the same.
#  2. Columns F ("化学纤维产销率") and G ("化学纤维销售量") are removed
#     entirely (they duplicated B/E except for the very first sub-row of
#     each year).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap the "B" and "C" sub-rows within every 4-row year block -------
# Year blocks start at row 2 and repeat every 4 rows through row 69.
# Within a block: first row = "A", second = "B", third = "C", fourth = "D".
# We swap the content of the "B" row and "C" row (columns A-E). For the
# blocks that also carry a value in column D (rows 54-69), D is swapped
# too; for the earlier blocks D is empty on both sides, so it is left
# untouched to avoid turning an empty-text cell into a blank cell.

for ($blockStart = 2; $blockStart -le 66; $blockStart += 4) {
    $rowB = $blockStart + 1
    $rowC = $blockStart + 2

    $rangeB = $ws.Range("A" + $rowB + ":E" + $rowB)
    $rangeC = $ws.Range("A" + $rowC + ":E" + $rowC)

    $valB = $rangeB.Value2
    $valC = $rangeC.Value2

    if ($blockStart -ge 54) {
        # Column D holds real data here too - swap the whole A:E range.
        $rangeB.Value2 = $valC
        $rangeC.Value2 = $valB
    }
    else {
        # Column D is blank for both rows in this block - swap only the
        # columns that actually carry data (A, B, C, E) and leave D as-is.
        $ws.Range("A" + $rowB).Value2 = $valC[1,1]
        $ws.Range("B" + $rowB).Value2 = $valC[1,2]
        $ws.Range("C" + $rowB).Value2 = $valC[1,3]
        $ws.Range("E" + $rowB).Value2 = $valC[1,5]

        $ws.Range("A" + $rowC).Value2 = $valB[1,1]
        $ws.Range("B" + $rowC).Value2 = $valB[1,2]
        $ws.Range("C" + $rowC).Value2 = $valB[1,3]
        $ws.Range("E" + $rowC).Value2 = $valB[1,5]
    }
}

# --- 2. Delete columns F and G (化学纤维产销率 / 化学纤维销售量) ----------
$ws.Range("F1:G69").EntireColumn.Delete()
